$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New row labels (B4:B31), post shared-string shift ----
$labelsText = @"
Holden
Rizzie Spiral
RotRing OmegaMax-90
Equal Angle
Tilt Rotate
CLR
Rizzie Hex
Matthies Hex
Tilt Rotate_Partial
RotRing OmegaMax-60
Equal Angle_Partial
Rizzie Hex_Partial
ND Single
RD Single
TD Single
Morris Single
Ring Perpendicular to ND
Ring Perpendicular to RD
Ring Perpendicular to TD
OffsetFTD
OffsetATD
OffsetF45
OffsetA45
OffsetFRD
OffsetARD
Gaussian Quadrature
Michael-CCHex
Michael-SNHex
"@
$labels = $labelsText -split "`n"

# ---- New numeric grid (C4:T31) ----
$dataText = @"
0.9999550003462915,1.000098998125962,0.9998800022410256,0.9999550003462915,0.9999550003462915,1.00017999838655,1.00017999838655,0.9999550003462915,0.9999550003462915,0.9999550003462915,1.000067499366421,1.000067499366421,1.000077998952934,1.000029999693044,1.000029999693044,1.000011249856356,1.000011249856356,1.000003999965402
1.000744684091735,0.9983617020121687,1.001985812962056,1.000744684091735,1.000744684091735,0.9970212724642875,0.9970212724642875,1.000744684091735,1.000744684091735,1.000744684091735,0.9988829782780114,0.9988829782780114,0.9987092195227305,0.9995035468825861,0.9995035468825861,0.9998138311848734,0.9998138311848734,0.9999338066189529
1.000253588620126,0.9994421055661601,1.000676239334453,1.000253588620126,1.000253588620126,0.9989856428938787,0.9989856428938787,1.000253588620126,1.000253588620126,1.000253588620126,0.9996196157570025,0.9996196157570025,0.9995604456933883,0.9998309400447104,0.9998309400447104,0.9999366021885644,0.9999366021885644,0.9999774589424785
0.999961371822766,1.000084979762251,0.9998969928962534,0.999961371822766,0.999961371822766,1.000154512384722,1.000154512384722,0.999961371822766,0.999961371822766,0.999961371822766,1.000057942103744,1.000057942103744,1.00006695465658,1.000025752010085,1.000025752010085,1.000009656963255,1.000009656963255,1.000003433418587
0.9998371828182898,1.000358175293824,0.9995658245880339,0.9998371828182898,0.9998371828182898,1.000651278723092,1.000651278723092,0.9998371828182898,0.9998371828182898,0.9998371828182898,1.000244230770691,1.000244230770691,1.000282212278402,1.000108548119891,1.000108548119891,1.00004070679449,1.00004070679449,1.000014471176637
0.9999979580850812,1.000004492020708,0.9999945589480972,0.9999979580850812,0.9999979580850812,1.000008162352322,1.000008162352322,0.9999979580850812,0.9999979580850812,0.9999979580850812,1.000003060218702,1.000003060218702,1.000003537486037,1.000001359507495,1.000001359507495,1.000000509151891,1.000000509151891,1.000000181262729
0.9999995939277415,1.000000894821559,0.9999989191671756,0.9999995939277415,0.9999995939277415,1.000001622991025,1.000001622991025,0.9999995939277415,0.9999995939277415,0.9999995939277415,1.000000608459383,1.000000608459383,1.000000703913442,1.000000270282169,1.000000270282169,1.000000101193562,1.000000101193562,1.000000036460498
0.9999910431493483,1.000019704778302,0.9999761180834331,0.9999910431493483,0.9999910431493483,1.000035823509013,1.000035823509013,0.9999910431493483,0.9999910431493483,0.9999910431493483,1.00001343332918,1.00001343332918,1.000015523812221,1.000005969935903,1.000005969935903,1.000002238239264,1.000002238239264,1.000000795969799
0.9998347407490412,1.000363548215615,0.9995593108408519,0.9998347407490412,0.9998347407490412,1.000661048732806,1.000661048732806,0.9998347407490412,0.9998347407490412,0.9998347407490412,1.000247894740924,1.000247894740924,1.000286445899154,1.00011017674363,1.00011017674363,1.000041317744983,1.000041317744983,1.0000146883394
1.000189303149439,0.9995835325477828,1.00050481415281,1.000189303149439,1.000189303149439,0.9992427821117049,0.9992427821117049,1.000189303149439,1.000189303149439,1.000189303149439,0.9997160426305721,0.9997160426305721,0.9996718726029757,0.9998737961368612,0.9998737961368612,0.9999526728900057,0.9999526728900057,0.999983173043436
0.9999578117157879,1.000092812305265,0.9998874965684197,0.9999578117157879,0.9999578117157879,1.000168756673681,1.000168756673681,0.9999578117157879,0.9999578117157879,0.9999578117157879,1.000063284194734,1.000063284194734,1.000073126898245,1.000028126701752,1.000028126701752,1.000010547955261,1.000010547955261,1.000003750115788
1.000032503494261,0.9999284978645003,1.000086672060078,1.000032503494261,1.000032503494261,0.9998699909465619,0.9998699909465619,1.000032503494261,1.000032503494261,1.000032503494261,0.9999512472204113,0.9999512472204113,0.9999436641017744,0.9999783326450279,0.9999783326450279,0.999991875357336,0.999991875357336,0.9999971118923203
0.9997144800000002,1.0006281,0.9992386200000002,0.9997144800000002,0.9997144800000002,1.0011421,1.0011421,0.9997144800000002,0.9997144800000002,0.9997144800000002,1.00042829,1.00042829,1.000494893333334,1.000190353333334,1.000190353333334,1.000071385,1.000071385,1.000025376666667
1.0010055,0.99778785,1.0026814,1.0010055,1.0010055,0.9959779000000001,0.9959779000000001,1.0010055,1.0010055,1.0010055,0.9984917,0.9984917,0.9982570833333333,0.9993296333333334,0.9993296333333334,0.9997486,0.9997486,0.9999106083333333
1.0018487,0.9959328,1.0049299,1.0018487,1.0018487,0.99260508,0.99260508,1.0018487,1.0018487,1.0018487,0.9972268900000001,0.9972268900000001,0.9967955266666667,0.9987674933333334,0.9987674933333334,0.9995377950000001,0.9995377950000001,0.9998356466666666
1.0007084,0.9984415900000001,1.001889,1.0007084,1.0007084,0.99716654,0.99716654,1.0007084,1.0007084,1.0007084,0.99893747,0.99893747,0.9987721766666667,0.9995277799999999,0.9995277799999999,0.9998229349999999,0.9998229349999999,0.9999370549999999
0.9999081668493154,1.000202032191781,0.9997551167123293,0.9999081668493154,0.9999081668493154,1.000367311506849,1.000367311506849,0.9999081668493154,0.9999081668493154,0.9999081668493154,1.000137739178082,1.000137739178082,1.000159170182648,1.000061215068493,1.000061215068493,1.000022953013699,1.000022953013699,1.000008160159817
1.000398363157895,0.9991235978947368,1.001062303157895,1.000398363157895,1.000398363157895,0.9984065526315791,0.9984065526315791,1.000398363157895,1.000398363157895,1.000398363157895,0.9994024578947369,0.9994024578947369,0.9993095045614035,0.9997344263157896,0.9997344263157896,0.9999004105263158,0.9999004105263158,0.9999645905263158
1.000675734210526,0.9985133868421052,1.001801949473684,1.000675734210526,1.000675734210526,0.9972970689473682,0.9972970689473682,1.000675734210526,1.000675734210526,1.000675734210526,0.9989864015789474,0.9989864015789474,0.9988287299999999,0.9995495124561403,0.9995495124561403,0.9998310678947369,0.9998310678947369,0.9999399346491228
0.9997681067625329,1.000510167563672,0.99938161847076,0.9997681067625329,0.9997681067625329,1.000927579456783,1.000927579456783,0.9997681067625329,0.9997681067625329,0.9997681067625329,1.000347843109658,1.000347843109658,1.000401951260996,1.000154597660616,1.000154597660616,1.000057974936095,1.000057974936095,1.000020614296469
0.9999367463468105,1.000139159407399,0.9998313272358839,0.9999367463468105,0.9999367463468105,1.000253014939221,1.000253014939221,0.9999367463468105,0.9999367463468105,0.9999367463468105,1.000094880643016,1.000094880643016,1.000109640231144,1.000042169210947,1.000042169210947,1.000015813494913,1.000015813494913,1.000005623437156
0.9999073222309677,1.00020388253831,0.9997528702554451,0.9999073222309677,0.9999073222309677,1.000370698243154,1.000370698243154,0.9999073222309677,0.9999073222309677,0.9999073222309677,1.000139010237061,1.000139010237061,1.000160634337477,1.000061780901696,1.000061780901696,1.000023166234014,1.000023166234014,1.000008236288302
0.999975776093317,1.000053302596141,0.999935393613998,0.999975776093317,0.999975776093317,1.000096907769418,1.000096907769418,0.999975776093317,0.999975776093317,0.999975776093317,1.000036341931367,1.000036341931367,1.000041995486292,1.000016153318684,1.000016153318684,1.000006059012342,1.000006059012342,1.000002155376585
1.000566909059617,0.9987528015024913,1.001511741381665,1.000566909059617,1.000566909059617,0.9977323649622599,0.9977323649622599,1.000566909059617,1.000566909059617,1.000566909059617,0.9991496370109387,0.9991496370109387,0.999017358508123,0.999622061027165,0.999622061027165,0.9998582730352781,0.9998582730352781,0.9999496058375449
1.000150352930528,0.9996692310591436,1.000400931423548,1.000150352930528,1.000150352930528,0.9993985986976708,0.9993985986976708,1.000150352930528,1.000150352930528,1.000150352930528,0.9997744758140993,0.9997744758140993,0.9997393942291141,0.9998997681862422,0.9998997681862422,0.9999624143723136,0.9999624143723136,0.9999866366619909
0.999982979861282,1.000037448706302,0.9999546326065073,0.999982979861282,0.999982979861282,1.00006806604681,1.00006806604681,0.999982979861282,0.999982979861282,0.999982979861282,1.000025522954046,1.000025522954046,1.000029498204798,1.000011341923124,1.000011341923124,1.000004251407664,1.000004251407664,1.000001514490578
1.000110357171643,0.9997572053302179,1.000294281822538,1.000110357171643,1.000110357171643,0.9995585754768366,0.9995585754768366,1.000110357171643,1.000110357171643,1.000110357171643,0.9998344663242398,0.9998344663242398,0.9998087126595658,0.9999264299400409,0.9999264299400409,0.9999724117479414,0.9999724117479414,0.9999901890240869
1.000157341691453,0.9996538208582806,1.000419586712324,1.000157341691453,1.000157341691453,0.9993706203664213,0.9993706203664213,1.000157341691453,1.000157341691453,1.000157341691453,0.9997639810289369,0.9997639810289369,0.9997272609720514,0.9998951012497755,0.9998951012497755,0.9999606613601948,0.9999606613601948,0.9999860088352306
"@
$dataLines = $dataText -split "`n"

$nrows = $dataLines.Count
$ncols = ($dataLines[0] -split ",").Count
$data = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
    $parts = $dataLines[$i] -split ","
    for ($j = 0; $j -lt $ncols; $j++) {
        $data[$i,$j] = [double]$parts[$j]
    }
}

# Write B4:B31 labels
for ($i = 0; $i -lt $labels.Count; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 2).Value2 = $labels[$i]
}

# Write C4:T31 numeric grid in one shot
$ws.Range("C4:T31").Value2 = $data

# Write A30/A31 row index values
$ws.Cells.Item(30, 1).Value2 = 28
$ws.Cells.Item(31, 1).Value2 = 29

# Match the bold/border/center style used by the other A-column header cells
$ws.Cells.Item(30, 1).Font.Bold = $true
$ws.Cells.Item(30, 1).HorizontalAlignment = -4108
$ws.Cells.Item(30, 1).VerticalAlignment = -4160
$ws.Cells.Item(30, 1).Borders.LineStyle = 1

$ws.Cells.Item(31, 1).Font.Bold = $true
$ws.Cells.Item(31, 1).HorizontalAlignment = -4108
$ws.Cells.Item(31, 1).VerticalAlignment = -4160
$ws.Cells.Item(31, 1).Borders.LineStyle = 1

Write-Output "Done"
